# Release 0.0.0 - Finish script
# Fill in the WhatsApp message table (rows 2-4) with sample data, keep the
# phone-number column formatted as Text, wrap the message column, refresh
# the comment that reminds to keep the phone column as Text, resize the
# columns to fit their content and restore the view (zoom/selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : User1 -----------------------------------------------------
$ws.Range("A2").Value = "User1"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "+51xxxxxxxxx"
$ws.Range("C2").Value = "text1"

# --- Row 3 : User2 -------------------------------------------------------
$ws.Range("A3").Value = "User2"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "+51xxxxxxxxx"
$ws.Range("C3").Value = "text2"

# --- Row 4 : User2 (second message) --------------------------------------
$ws.Range("A4").Value = "User2"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "+51xxxxxxxxx"
$ws.Range("C4").Value = "text3"

# --- Image column, filled in afterwards for all three rows ---------------
$ws.Range("D2").Value = "robot_hi.png"
$ws.Range("D3").Value = "robot_hi.png"
$ws.Range("D4").Value = "robot_hi.png"

# Keep the rest of the "Celular" column as Text as well (rows with no data yet)
$ws.Range("B5:B10").NumberFormat = "@"

# Wrap the message text so long messages are readable
$ws.Range("C2:C4").WrapText = $true

# Remind to keep the phone number column as Text
$cmt = $ws.Range("B1").AddComment("Sergio Gutiérrez Sanchez:" + [char]10 + "Mantener en Tipo Texto")

# Resize the columns to fit their new content
$ws.Columns("A:D").AutoFit() | Out-Null

# Restore view: zoom out a bit and move the selection
$excel.ActiveWindow.Zoom = 130
$ws.Range("H10").Select() | Out-Null
